$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The existing row 4 (Vega Modelo de Temuco / Breva, dated 2023-11-09 with
# Volumen=25 and Origen="Provincia de San Felipe de Aconcagua") becomes a new
# weekly report row 5, while row 4 itself is refreshed with this week's data
# (new date, new Volumen, new Origen).

# 1) Copy the current row 4 (last data row) down into the new row 5, keeping
#    values and the date number format intact.
for ($col = 1; $col -le 20; $col++) {
    $src = $ws.Cells.Item(4, $col)
    $dst = $ws.Cells.Item(5, $col)
    $dst.Value2 = $src.Value2
    if ($col -eq 4) {
        $dst.NumberFormat = $src.NumberFormat
    }
}

# 2) Overwrite row 4 with the new week's values.
$ws.Cells.Item(4, 4).Value2 = 45244          # D4 - Fecha
$ws.Cells.Item(4, 13).Value2 = 70            # M4 - Volumen
$ws.Cells.Item(4, 18).Value2 = "Región Metropolitana"  # R4 - Origen
